$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking crypto snapshot refresh: column D is Price, column E is the 1-hour
# percent change. A handful of Price values are purely numeric-looking text (e.g.
# "0.9997", "15.10") that Excel would otherwise silently convert to a number,
# dropping formatting such as trailing zeros. A leading apostrophe forces those to
# stay text (same as typing `'0.9997` by hand); resetting the style afterwards
# clears the transient "quote prefix" text-entry formatting it leaves behind.
$apostrophe = [string][char]39

$ws.Range("D2").Value = '30.315.89'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.934.73'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("D4").Value = $apostrophe + '0.9997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.77%  '
$ws.Range("D5").Value = $apostrophe + '250.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("D6").Value = $apostrophe + '0.7258'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.97%  '
$ws.Range("D7").Value = $apostrophe + '0.9992'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.76%  '
$ws.Range("D8").Value = $apostrophe + '0.3307'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.52%  '
$ws.Range("D9").Value = $apostrophe + '27.95'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.99%  '
$ws.Range("D10").Value = $apostrophe + '0.07280'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.22%  '
$ws.Range("D11").Value = $apostrophe + '0.8079'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.85%  '
$ws.Range("D12").Value = $apostrophe + '0.08100'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.32%  '
$ws.Range("D13").Value = '1.935.70'
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("D14").Value = $apostrophe + '5.481'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.84%  '
$ws.Range("D15").Value = $apostrophe + '94.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = $apostrophe + '15.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.62%  '
$ws.Range("D17").Value = '30.318.16'
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").Value = $apostrophe + '254.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.34%  '
$ws.Range("D19").Value = $apostrophe + '0.000008230'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.54%  '
$ws.Range("D20").Value = $apostrophe + '5.804'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("D21").Value = '2.187.86'
$ws.Range("E21").Value = '  -0.57%  '
$ws.Range("D22").Value = $apostrophe + '0.9990'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.70%  '
$ws.Range("D23").Value = $apostrophe + '0.9999'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = $apostrophe + '6.965'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.25%  '
$ws.Range("D25").Value = $apostrophe + '9.766'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("D26").Value = $apostrophe + '165.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.14%  '
$ws.Range("D27").Value = $apostrophe + '2.358'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.13%  '
$ws.Range("E28").Value = '  +3.08%  '
$ws.Range("D29").Value = $apostrophe + '0.1308'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").Value = $apostrophe + '1.539'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("D32").Value = $apostrophe + '4.436'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.74%  '
$ws.Range("D33").Value = $apostrophe + '4.200'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("D34").Value = $apostrophe + '0.05264'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.71%  '
$ws.Range("D35").Value = $apostrophe + '1.272'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.11%  '
$ws.Range("D36").Value = $apostrophe + '0.7498'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("D37").Value = $apostrophe + '2.766'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.34%  '
$ws.Range("E38").Value = '  +2.58%  '
$ws.Range("D39").Value = $apostrophe + '2.809'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.90%  '
$ws.Range("D40").Value = $apostrophe + '79.26'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D41").Value = $apostrophe + '6.444'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.17%  '
$ws.Range("D42").Value = $apostrophe + '0.4539'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.81%  '
$ws.Range("D43").Value = $apostrophe + '2.033'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("D44").Value = $apostrophe + '0.8451'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("D45").Value = $apostrophe + '0.9995'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("D46").Value = $apostrophe + '101.76'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("D47").Value = $apostrophe + '9.749'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("E48").Value = '  +2.50%  '
$ws.Range("D49").Value = $apostrophe + '36.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.45%  '
$ws.Range("D50").Value = $apostrophe + '0.4193'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.93%  '
$ws.Range("D51").Value = $apostrophe + '0.06045'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.68%  '
